$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 92, pushing existing rows 92..134 down to 93..135
$ws.Rows.Item(92).Insert()

# Populate the newly inserted row 92 with the new record
$ws.Cells.Item(92, 1).Value = 1
$ws.Cells.Item(92, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(92, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(92, 4).Value = 45009
$ws.Cells.Item(92, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(92, 5).Value = 15
$ws.Cells.Item(92, 6).Value = "Fruta"
$ws.Cells.Item(92, 7).Value = 100102
$ws.Cells.Item(92, 8).Value = "Cítricos"
$ws.Cells.Item(92, 9).Value = 100102005
$ws.Cells.Item(92, 10).Value = "Naranja"
$ws.Cells.Item(92, 11).Value = "Valencia"
$ws.Cells.Item(92, 12).Value = "Tercera"
$ws.Cells.Item(92, 13).Value = 270
$ws.Cells.Item(92, 14).Value = 1100
$ws.Cells.Item(92, 15).Value = 1200
$ws.Cells.Item(92, 16).Value = 1150
$ws.Cells.Item(92, 17).Value = "`$/kilo (en caja de 20 kilos)"
$ws.Cells.Item(92, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(92, 19).Value = 1150
$ws.Cells.Item(92, 20).Value = 1
